# Generate Report for Handoff
# The file "f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.md" has been handed off again:
#  - its Status becomes "Ready for handoff" on the Overview sheet and on each
#    language sheet (zh-cn, de-de)
#  - the "Latest Handoff Datetime" column is refreshed with the new handoff
#    timestamp for each language sheet (both data rows in that column share
#    the same stamp, matching the source report)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$readyStatus = "Ready for handoff"

# Overview sheet: row 3 corresponds to f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.md
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus

# zh-cn sheet: row 3 corresponds to f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.md
$wsZhCn.Range("B3").Value = $readyStatus
$wsZhCn.Range("D2").Value = "2016-03-09 16:13:06"
$wsZhCn.Range("D3").Value = "2016-03-09 16:13:06"

# de-de sheet: row 3 corresponds to f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.md
$wsDeDe.Range("B3").Value = $readyStatus
$wsDeDe.Range("D2").Value = "2016-03-09 16:13:16"
$wsDeDe.Range("D3").Value = "2016-03-09 16:13:16"
